# Refresh the cryptocurrency price/volume table (cols D & E, rows 2-51)
# with the latest scraped values, matching the source GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text for each changed cell. Values are entered as literal text, exactly
# as the scraper writes them (thousands separated by "." for some coins,
# trailing zeros kept in places, "%" rows padded with two spaces each side).
$newValues = [ordered]@{
    "D2" = "29.241.79"
    "E2" = "  +0.32%  "
    "D3" = "1.855.95"
    "E3" = "  +0.21%  "
    "E4" = "  +0.13%  "
    "D5" = "0.7022"
    "E5" = "  +2.81%  "
    "D6" = "238.03"
    "E6" = "  +0.28%  "
    "E7" = "  +0.00%  "
    "D8" = "0.07996"
    "E8" = "  +4.67%  "
    "D9" = "0.3029"
    "E9" = "  -0.31%  "
    "D10" = "23.62"
    "E10" = "  +2.27%  "
    "D11" = "0.08186"
    "E11" = "  +0.74%  "
    "D12" = "1.847.46"
    "E12" = "  +0.22%  "
    "D13" = "5.200"
    "D14" = "0.7071"
    "E14" = "  -2.02%  "
    "D15" = "89.78"
    "E15" = "  +0.50%  "
    "D16" = "29.163.48"
    "E16" = "  +0.04%  "
    "D17" = "5.826"
    "E17" = "  +2.04%  "
    "D18" = "0.000007855"
    "E18" = "  +0.88%  "
    "D19" = "13.23"
    "E19" = "  +0.31%  "
    "D20" = "236.72"
    "E20" = "  +1.62%  "
    "D21" = "0.9997"
    "E21" = "  -0.10%  "
    "D22" = "1.001"
    "E22" = "  +0.02%  "
    "D23" = "2.071.75"
    "E23" = "  -1.57%  "
    "D24" = "7.525"
    "E24" = "  +1.62%  "
    "D25" = "163.26"
    "E25" = "  +1.08%  "
    "D26" = "8.890"
    "E26" = "  -0.45%  "
    "D27" = "0.1409"
    "E27" = "  -0.80%  "
    "D28" = "18.12"
    "E28" = "  +0.77%  "
    "D29" = "1.909"
    "E29" = "  -2.11%  "
    "D30" = "1.402"
    "E30" = "  +0.37%  "
    "D31" = "1.472"
    "E31" = "  -0.60%  "
    "D32" = "4.354"
    "E32" = "  -3.22%  "
    "D33" = "4.029"
    "E33" = "  +0.68%  "
    "D34" = "0.05191"
    "E34" = "  +0.92%  "
    "D35" = "1.166"
    "E35" = "  -0.96%  "
    "D36" = "0.7151"
    "E36" = "  +1.85%  "
    "D37" = "0.9977"
    "E37" = "  -2.45%  "
    "D38" = "2.680"
    "E38" = "  +0.37%  "
    "E39" = "  +0.25%  "
    "D40" = "2.721"
    "E40" = "  +1.67%  "
    "D41" = "0.9356"
    "E41" = "  +3.46%  "
    "D42" = "1.150.94"
    "E42" = "  +4.37%  "
    "D43" = "5.997"
    "E43" = "  +0.63%  "
    "E44" = "  -0.21%  "
    "D45" = "70.20"
    "E45" = "  +0.69%  "
    "D46" = "1.001"
    "E46" = "  +0.03%  "
    "D47" = "102.92"
    "E47" = "  +0.59%  "
    "D48" = "0.5290"
    "E48" = "  -3.21%  "
    "D49" = "1.745"
    "E49" = "  -1.45%  "
    "D50" = "1.999.94"
    "E50" = "  +0.19%  "
    "D51" = "9.152"
    "E51" = "  +0.33%  "
}

# Cells whose new text is a plain number (e.g. "238.03") would otherwise be
# auto-recognised by Excel as a Number and lose their exact text representation
# (e.g. "5.200" -> 5.2). Force those specific cells to Text format first so the
# literal string is preserved, matching how the sheet was originally authored.
$forceTextCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D49", "D51"
)
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = $newValues[$ref]
}

